$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '25.991.49'
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '1.629.55'
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.78%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '214.01'
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '0.504'
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.72%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.250'
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.95%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.0619'
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.04%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '18.48'
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.52%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.0787'
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.02%  '
$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = '@'
$c.Value = 'WrappedEther'
$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '1.750.83'
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = '@'
$c.Value = '  +4.14%  '
$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = '@'
$c.Value = 'WrappedliquidstakedEther2.0'
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '1.854.33'
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '4.19'
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.98%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.528'
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.93%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '25.983.28'
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.28%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0739'
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.21%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '61.33'
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.25%  '
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '192.27'
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.08%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '4.24'
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.66%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '9.60'
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.33%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '6.08'
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.93%  '
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.76%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '143.56'
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.37%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '1.74'
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.98%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '6.73'
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.06%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '15.22'
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.0484'
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.02%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '3.13'
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = '@'
$c.Value = '  -4.01%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '3.12'
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.30%  '
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.45%  '
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.82%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '1.127.69'
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.33%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.854'
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.65%  '
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = '@'
$c.Value = '  -1.18%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.522'
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.40%  '
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = '@'
$c.Value = '  -2.26%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '98.20'
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.86%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '1.764.08'
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '0.766'
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.94%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '5.18'
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = '@'
$c.Value = '  -5.25%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₆0104'
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = '@'
$c.Value = '  -11.02%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.0532'
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = '@'
$c.Value = '  +1.87%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '54.46'
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.61%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '1.48'
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.95%  '
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = '@'
$c.Value = '  -0.31%  '
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = '@'
$c.Value = '  +0.18%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '7.47'
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = '@'
$c.Value = '  -3.76%  '
